# Updates the "Estado de Cuenta" worksheet with a refreshed database extract:
# - Worker #2 (row 17) changes from WILLIAMS JOSE GUEVARA GOMEZ (73200363) to
#   AURY LUZ GAVIRIA PUERTA (1048434925), whose periods now run 2404..2508
#   (ascending) instead of the previous 2404..2507 block (which ran
#   descending from row 18 downward).
# - Summary cells (Valor Mora total, worker count, period count) are
#   refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary block -------------------------------------------------------
$ws.Range("E11").Value = 917479      # VALOR MORA total
$ws.Range("C13").Value = 2           # Cant. Trabajadores
$ws.Range("F13").Value = 18          # Cant. Periodos

# --- Data rows 17-33: second worker's full period history ---------------
# Row 16 (first worker, ALVARO DE JESUS MALLARINO SANCHEZ) is unchanged.

$docType = "CC"
$docNum = "1048434925"
$name = "AURY LUZ GAVIRIA PUERTA"
$salario = 1300000

$periods = @(
    @{ Row = 17; Periodo = "2404"; Mora = 27733 },
    @{ Row = 18; Periodo = "2405"; Mora = 52000 },
    @{ Row = 19; Periodo = "2406"; Mora = 52000 },
    @{ Row = 20; Periodo = "2407"; Mora = 52000 },
    @{ Row = 21; Periodo = "2408"; Mora = 52000 },
    @{ Row = 22; Periodo = "2409"; Mora = 52000 },
    @{ Row = 23; Periodo = "2410"; Mora = 52000 },
    @{ Row = 24; Periodo = "2411"; Mora = 52000 },
    @{ Row = 25; Periodo = "2412"; Mora = 52000 },
    @{ Row = 26; Periodo = "2501"; Mora = 52000 },
    @{ Row = 27; Periodo = "2502"; Mora = 52000 },
    @{ Row = 28; Periodo = "2503"; Mora = 52000 },
    @{ Row = 29; Periodo = "2504"; Mora = 52000 },
    @{ Row = 30; Periodo = "2505"; Mora = 52000 },
    @{ Row = 31; Periodo = "2506"; Mora = 52000 },
    @{ Row = 32; Periodo = "2507"; Mora = 52000 },
    @{ Row = 33; Periodo = "2508"; Mora = 52000 }
)

foreach ($p in $periods) {
    $r = $p.Row
    $ws.Range("B$r").Value = $docType
    $ws.Range("C$r").Value = $docNum
    $ws.Range("D$r").Value = $name
    $ws.Range("E$r").Value = $p.Periodo
    $ws.Range("F$r").Value = $p.Mora
    $ws.Range("G$r").Value = $salario
}
